# Apply edits to Sheet1: rename a header, add a Break Time column,
# renumber/merge the pomodoro rows and replace start/end time values
# with plain text time-of-day strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two pomodoro-session rows (old rows 10 and 11) were consolidated into one
# row, so delete a row to shrink the sheet from 22 rows to 21 rows.
$ws.Rows.Item(11).Delete()

# Header row updates
$ws.Range("C1").Value = "Minutes"
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Break Time"

# Row 2
$ws.Range("A2").Value = "Introduction"
$ws.Range("B2").Value = 44572.54166666666
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "13:00:00"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "13:30:00"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = "13:30:00"
$ws.Range("H2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = "Kubernetes Overview"
$ws.Range("B3").Value = 44572.54166666666
$ws.Range("C3").Value = 21
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "13:00:00"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "13:30:00"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = "13:30:00"
$ws.Range("H3").Style = "Normal"

# Row 4
$ws.Range("A4").Value = "Kubernetes Overview"
$ws.Range("B4").Value = 44573.54166666666
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 31
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = "13:00:00"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "13:30:00"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "13:30:00"
$ws.Range("H4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = "Setup Kubernetes"
$ws.Range("B5").Value = 44573.54166666666
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 51
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "13:00:00"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "13:30:00"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = "13:30:00"
$ws.Range("H5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = "Kubernetes Concepts"
$ws.Range("B6").Value = 44573.54166666666
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "13:00:00"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "13:30:00"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = "13:30:00"
$ws.Range("H6").Style = "Normal"

# Row 7
$ws.Range("A7").Value = "Kubernetes Concepts"
$ws.Range("B7").Value = 44574.54166666666
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 64
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = "13:00:00"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "13:30:00"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "13:30:00"
$ws.Range("H7").Style = "Normal"

# Row 8
$ws.Range("A8").Value = "YAML Introduction"
$ws.Range("B8").Value = 44574.54166666666
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 72
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = "13:00:00"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "13:30:00"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = "13:30:00"
$ws.Range("H8").Style = "Normal"

# Row 9
$ws.Range("A9").Value = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Range("B9").Value = 44574.54166666666
$ws.Range("C9").Value = 18
$ws.Range("D9").Value = 90
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "13:00:00"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "13:30:00"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = "13:30:00"
$ws.Range("H9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Range("B10").Value = 44575.54166666666
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 120
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "13:00:00"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "13:30:00"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = "13:30:00"
$ws.Range("H10").Style = "Normal"

# Row 11
$ws.Range("A11").Value = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Range("B11").Value = 44576.54166666666
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 150
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "13:00:00"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "13:30:00"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = "13:30:00"
$ws.Range("H11").Style = "Normal"

# Row 12
$ws.Range("A12").Value = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Range("B12").Value = 44577.54166666666
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 180
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = "13:00:00"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "13:30:00"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = "13:30:00"
$ws.Range("H12").Style = "Normal"

# Row 13
$ws.Range("A13").Value = "Kubernets Concepts - PODs, ReplicaSets, Deployments"
$ws.Range("B13").Value = 44578.54166666666
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 200
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "13:00:00"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = "13:30:00"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = "13:30:00"
$ws.Range("H13").Style = "Normal"

# Row 14
$ws.Range("A14").Value = "Networking in Kubernetes"
$ws.Range("B14").Value = 44578.54166666666
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 205
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = "13:00:00"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = "13:30:00"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = "13:30:00"
$ws.Range("H14").Style = "Normal"

# Row 15
$ws.Range("A15").Value = "Services"
$ws.Range("B15").Value = 44578.54166666666
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 210
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = "13:00:00"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = "13:30:00"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = "13:30:00"
$ws.Range("H15").Style = "Normal"

# Row 16
$ws.Range("A16").Value = "Services"
$ws.Range("B16").Value = 44579.54166666666
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 229
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = "13:00:00"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = "13:30:00"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = "13:30:00"
$ws.Range("H16").Style = "Normal"

# Row 17
$ws.Range("A17").Value = "Microservices Architechture"
$ws.Range("B17").Value = 44579.54166666666
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 240
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = "13:00:00"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = "13:30:00"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = "13:30:00"
$ws.Range("H17").Style = "Normal"

# Row 18
$ws.Range("A18").Value = "Microservices Architechture"
$ws.Range("B18").Value = 44580.54166666666
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 270
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = "13:00:00"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = "13:30:00"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = "13:30:00"
$ws.Range("H18").Style = "Normal"

# Row 19
$ws.Range("A19").Value = "Microservices Architechture"
$ws.Range("B19").Value = 44581.54166666666
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 274
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = "13:00:00"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = "13:30:00"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value = "13:30:00"
$ws.Range("H19").Style = "Normal"

# Row 20
$ws.Range("A20").Value = "Kubernetes on the Cloud"
$ws.Range("B20").Value = 44581.54166666666
$ws.Range("C20").Value = 26
$ws.Range("D20").Value = 300
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = "13:00:00"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = "13:30:00"
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value = "13:30:00"
$ws.Range("H20").Style = "Normal"

# Row 21
$ws.Range("A21").Value = "Conclusion"
$ws.Range("B21").Value = 44582.54166666666
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 302
$ws.Range("E21").Value = 11
$ws.Range("F21").Value = "13:00:00"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = "13:30:00"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = "13:30:00"
$ws.Range("H21").Style = "Normal"

